$d = $word.ActiveDocument

# 1. "Giao tác ... tham chiếu đến TestCase02_T1.sql." -> "...TestCase03_T1.sql."
$p4 = $d.Paragraphs(4)
$ok1 = $p4.Range.Find.Execute(
    "TestCase02_T1.sql", $true, $false, $false, $false, $false,
    $true, 1, $false, "TestCase03_T1.sql", 2)
Write-Host "1 TestCase03_T1: $ok1"

# 2. "Giao tác ... tham chiếu đến TestCase02_T2.sql" -> "...TestCase03_T2.sql"
$p5 = $d.Paragraphs(5)
$ok2 = $p5.Range.Find.Execute(
    "TestCase02_T2.sql", $true, $false, $false, $false, $false,
    $true, 1, $false, "TestCase03_T2.sql", 2)
Write-Host "2 TestCase03_T2: $ok2"

# 3. T1 reads and now also SAVES the due date; introduces the course end-date sentence.
$p8 = $d.Paragraphs(8)
$ok3 = $p8.Range.Find.Execute(
    "Giao tác T1 đọc thời gian nộp của đồ án và in ra màn hình bằng lệnh PRINT. Thời gian nộp lúc này là 30/05/2012.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Giao tác T1 đọc thời gian nộp của đồ án và lưu lại. Giả sử thời gian nộp lúc này là 30/07/2012, ngày kết thúc môn học ứng với đồ án trên là 01/08/2012.",
    2)
Write-Host "3 PRINT paragraph: $ok3"

# 4. T1 checks validity "để" insert data (instead of ". Nếu hợp lệ thì")
$p9 = $d.Paragraphs(9)
$ok4 = $p9.Range.Find.Execute(
    "Giao tác T1 kiểm tra thời gian nộp đồ án xem có hợp lệ không. Nếu hợp lệ thì thêm dữ liệu vào bảng DE_SINHVIEN.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Giao tác T1 kiểm tra thời gian nộp đồ án xem có hợp lệ không để thêm dữ liệu vào bảng DE_SINHVIEN.",
    2)
Write-Host "4 kiem tra paragraph: $ok4"

# 5. T2 changes the due date to 29/07/2012 (was 29/05/2012)
$p10 = $d.Paragraphs(10)
$ok5 = $p10.Range.Find.Execute(
    "29/05/2012", $true, $false, $false, $false, $false,
    $true, 1, $false, "29/07/2012", 2)
Write-Host "5 T2 change date: $ok5"

# 6. Insert a new bullet paragraph before "Giao tác T1 đọc lại..." describing what T1 sees
#    after it re-reads/saves and prints the value it stored earlier.
$p11 = $d.Paragraphs(11)
$newRange = $p11.Range.Duplicate
$newRange.Collapse(1)
$newRange.InsertParagraphBefore()
$newPara = $d.Paragraphs(11)
$newPara.Range.Text = "Sau đó T1 thấy dữ liệu hợp lệ và xuất ra màn hình thời hạn nộp đã lưu trước đó (bằng lệnh PRINT) là 30/07/2012."
Write-Host "6 new paragraph inserted, count=$($d.Paragraphs.Count)"

# 7. The re-read value becomes 29/07/2012 (was 29/05/2012); paragraph index shifted by +1 to 12.
$p12 = $d.Paragraphs(12)
$ok7 = $p12.Range.Find.Execute(
    "29/05/2012", $true, $false, $false, $false, $false,
    $true, 1, $false, "29/07/2012", 2)
Write-Host "7 T1 re-read date: $ok7"

# 8. Kết luận paragraph - clarify "(trước và sau khi T2 chạy)" and "giá trị" instead of "dữ liệu".
$p13 = $d.Paragraphs(13)
$ok8 = $p13.Range.Find.Execute(
    "thời gian nộp do giao tác T1 đọc ở hai lần có dữ liệu khác nhau.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "thời gian nộp do giao tác T1 đọc ở hai lần (trước và sau khi T2 chạy) có giá trị khác nhau.",
    2)
Write-Host "8 Ket luan paragraph: $ok8"

# 9. Cách khắc phục paragraph - append explicit resulting/locked-update dates to the conclusion.
$p14 = $d.Paragraphs(14)
$ok9 = $p14.Range.Find.Execute(
    "sau đó mới thực hiện việc cập nhật thời gian nộp. Do vậy khi T1 xuất kết quả thời gian nộp ở hai lần sẽ ra giống nhau.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "sau đó mới thực hiện việc cập nhật thời gian nộp. Do vậy khi T1 xuất kết quả thời gian nộp ở hai lần sẽ ra giống nhau là 30/07/2012. Sau khi T1 kết thúc, T2 mới được phép cập nhật thời gian nộp thành 29/07/2012.",
    2)
Write-Host "9 Cach khac phuc paragraph: $ok9"
